# AP: Problemas ao crypto_excel
#
# Reworks the data grid on the active sheet:
#  - Row 1: shifts the label/index values over by one column, introducing a
#    new A1 text label and extending the numeric sequence in F1:K1 through 5.
#  - Rows 2-20: re-keys the "teste/eu/sou/o/jonas/silva/teste/teste/aeeeo/
#    acabouaqui" sentence across A:K (column A had held a running numeric
#    index and loses both that number and its header-style formatting), and
#    appends a trailing "acaba" value in column K.
#  - Finally, turns on worksheet protection (matches Excel's default
#    "Protect Sheet" dialog options) with a password.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1: numeric cells first (F1:K1) so K1 still carries the original
# header style (s=1) to copy from when we fix up the text cells afterwards.
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = 0
$ws.Range("G1").Value = 1
$ws.Range("H1").Value = 2
$ws.Range("I1").Value = 3
$ws.Range("J1").Value = 4
$ws.Range("K1").Value = 5

# Text cells A1:E1 -- the leading apostrophe forces these numeric-looking
# strings to stay text; pasting K1's format back on top keeps the header
# style (border + bold + centered) instead of the transient quote-prefix
# style that assigning a quoted value creates.
$ws.Range("A1").Value = "'1"
$ws.Range("B1").Value = "'2"
$ws.Range("C1").Value = "'3"
$ws.Range("D1").Value = "'4"
$ws.Range("E1").Value = "'5"

$ws.Range("K1").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# Rows 2-20: same sentence in every row, column A drops its old numeric
# index + header style, column K gets a new trailing word.
# ---------------------------------------------------------------------------
$words = @("teste", "eu", "sou", "o", "jonas", "silva", "teste", "teste", "aeeeo", "acabouaqui", "acaba")
$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

for ($row = 2; $row -le 20; $row++) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $words[$i]
    }

    # Column A used to carry the bordered/bold/centered header style (s=1)
    # alongside its numeric index; strip that back to the plain default
    # style shared by the rest of the row (e.g. column B).
    $ws.Range("B$row").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Protect the sheet (default allow-list, matching the Excel "Protect Sheet"
# dialog with nothing toggled) behind a password.
# ---------------------------------------------------------------------------
$ws.Protect("CC3D")
